# "Generate Report for Handback" — records a failed handback validation
# (stale handback file) for 0ba0874e-67bf-40bf-8407-811cb0beebf2.md on both
# the zh-cn and de-de localization-status sheets.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7181042e34a780be450babe9e5c4cac7fc8dd569/e2e/0ba0874e-67bf-40bf-8407-811cb0beebf2.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dc5337ccdb700cb9a635f408bb6bdef1a1831d29/e2e/0ba0874e-67bf-40bf-8407-811cb0beebf2.md."
$latestTargetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dc5337ccdb700cb9a635f408bb6bdef1a1831d29/e2e/0ba0874e-67bf-40bf-8407-811cb0beebf2.md"
$latestTargetDisplay = "0ba0874e-67bf-40bf-8407-811cb0beebf2.md"

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")

# Widen the "Error Detail" column (P / col 16) so the long message is readable.
$wsZh.Columns.Item(16).ColumnWidth = 39.1671875

# Row 6 is 0ba0874e-67bf-40bf-8407-811cb0beebf2.md — stamp its "Latest Target
# File" with a hyperlink to the current/expected commit, fill in the matching
# handback xliff + timestamp, and record the validation error.
$wsZh.Hyperlinks.Add($wsZh.Range("I6"), $latestTargetUrl, [Type]::Missing, [Type]::Missing, $latestTargetDisplay)
$wsZh.Range("J6").Value = "0ba0874e-67bf-40bf-8407-811cb0beebf2.05886c9f81c78f4e0b1738b22e8a6ac80dd65499.zh-cn.xlf"
$wsZh.Range("K6").Value = "2016-08-16 06:38:30"
$wsZh.Range("P6").Value = $errorDetail

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Columns.Item(16).ColumnWidth = 39.1671875

$wsDe.Hyperlinks.Add($wsDe.Range("I6"), $latestTargetUrl, [Type]::Missing, [Type]::Missing, $latestTargetDisplay)
$wsDe.Range("J6").Value = "0ba0874e-67bf-40bf-8407-811cb0beebf2.05886c9f81c78f4e0b1738b22e8a6ac80dd65499.de-de.xlf"
$wsDe.Range("K6").Value = "2016-08-16 06:38:37"
$wsDe.Range("P6").Value = $errorDetail
